# Remove the trailing "Requisitos" heading paragraph and the following
# "LOQ4044 - ... (Requisito fraco)" bullet paragraph from the end of the
# document (the course no longer lists a prerequisite section).

$d = $word.ActiveDocument

$reqHeading = $null
$reqBullet = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.TrimEnd() -eq "Requisitos") {
        $reqHeading = $p
    } elseif ($t -like "LOQ4044*") {
        $reqBullet = $p
    }
}

if ($reqHeading -ne $null -and $reqBullet -ne $null) {
    $rng = $d.Range($reqHeading.Range.Start, $reqBullet.Range.End)
    $rng.Delete()
}
